# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, centered, bordered) from the existing
# "Unnamed: 28" header cell (AC1) onto the three new header cells so they
# match the look of the rest of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill the team record for every player row (2-41) with the 1991 Brewers
# record: 83 wins, 79 losses, 0 ties.
$ws.Range("AD2:AD41").Value = 83
$ws.Range("AE2:AE41").Value = 79
$ws.Range("AF2:AF41").Value = 0
